# Implemented and tested gamma functions.
# Updates the lgamma (row 34) and tgamma (row 35) benchmark rows with new
# measured values, and flags the two "ULP==0" sample-size cells (H34, H35)
# with the built-in "Bad" cell style (as Excel's conditional red/pink
# highlight) while keeping their original box borders intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Updated measured values - row 34 ("lgamma")
# ---------------------------------------------------------------------
$ws.Range("B34").Value = 17
$ws.Range("C34").Value = 16
$ws.Range("D34").Value = 16
$ws.Range("E34").Value = 16
$ws.Range("F34").Value = 31192
$ws.Range("G34").Value = 32155
$ws.Range("M34").Value = 9089
$ws.Range("N34").Value = 5915
$ws.Range("W34").Value = 10178
$ws.Range("X34").Value = 7494

# ---------------------------------------------------------------------
# 2. Updated measured values - row 35 ("tgamma")
# ---------------------------------------------------------------------
$ws.Range("B35").Value = 121
$ws.Range("C35").Value = 58
$ws.Range("D35").Value = 57
$ws.Range("E35").Value = 58
$ws.Range("F35").Value = 24722
$ws.Range("G35").Value = 24292
$ws.Range("M35").Value = 9583
$ws.Range("N35").Value = 5250
$ws.Range("W35").Value = 23479
$ws.Range("X35").Value = 15897

# ---------------------------------------------------------------------
# 3. Highlight H34 and H35 with the built-in "Bad" cell style, preserving
#    each cell's existing border.
# ---------------------------------------------------------------------
function Set-BadStyleKeepBorder($cell) {
    $edges = @(7, 8, 9, 10)
    $saved = @{}
    foreach ($e in $edges) {
        $b = $cell.Borders.Item($e)
        $saved[$e] = @{ LineStyle = $b.LineStyle(); Weight = $b.Weight(); ColorIndex = $b.ColorIndex() }
    }

    $cell.Style = "Bad"

    foreach ($e in $edges) {
        $s = $saved[$e]
        if ($s.LineStyle -ne -4142) {
            $b = $cell.Borders.Item($e)
            $b.LineStyle = $s.LineStyle
            $b.Weight = $s.Weight
            $b.ColorIndex = $s.ColorIndex
        }
    }
}

Set-BadStyleKeepBorder $ws.Range("H34")
Set-BadStyleKeepBorder $ws.Range("H35")

$wb.Save()
